$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 308369.03
$ws.Range("J17").Value = 308369.03
$ws.Range("L17").Value = 925107.0900000001
$ws.Range("N17").Value = -925443.0900000001
$ws.Range("H129").Value = 2058679
$ws.Range("J129").Value = 2179724.8
$ws.Range("L129").Value = 6539174.399999999
$ws.Range("N129").Value = -6549174.399999999
$ws.Range("H137").Value = 7938122.5
$ws.Range("I137").Value = 1354.9318
$ws.Range("J137").Value = 26318004
$ws.Range("K137").Value = 4064.7954
$ws.Range("L137").Value = 78954012
$ws.Range("M137").Value = -1514.7954
$ws.Range("N137").Value = -78959112

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1154.6285
$ws.Range("J74").Value = 919.3889
$ws.Range("L74").Value = 919.3889
$ws.Range("N74").Value = -2667.3889
$ws.Range("H77").Value = 1154.6285
$ws.Range("J77").Value = 919.3889
$ws.Range("L77").Value = 4596.944500000001
$ws.Range("N77").Value = -13332.9445
$ws.Range("H88").Value = 2538.2104
$ws.Range("I88").Value = 2583.0625
$ws.Range("J88").Value = 2299
$ws.Range("K88").Value = 2583.0625
$ws.Range("L88").Value = 2299
$ws.Range("M88").Value = -2177.0625
$ws.Range("N88").Value = -3111
$ws.Range("H91").Value = 2538.2104
$ws.Range("I91").Value = 2583.0625
$ws.Range("J91").Value = 2299
$ws.Range("K91").Value = 2583.0625
$ws.Range("L91").Value = 2299
$ws.Range("M91").Value = -1179.0625
$ws.Range("N91").Value = -5107
$ws.Range("H122").Value = 1170.6666
$ws.Range("I122").Value = 1170.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3511.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1061.9998
$ws.Range("N122").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2413.0833
$ws.Range("I86").Value = 2265.0715
$ws.Range("K86").Value = 2265.0715
$ws.Range("M86").Value = -1142.0715
$ws.Range("H89").Value = 2413.0833
$ws.Range("I89").Value = 2265.0715
$ws.Range("K89").Value = 11325.3575
$ws.Range("M89").Value = -5709.3575
$ws.Range("H99").Value = 58825760
$ws.Range("I99").Value = 83335470
$ws.Range("J99").Value = 2456.2
$ws.Range("K99").Value = 83335470
$ws.Range("L99").Value = 2456.2
$ws.Range("M99").Value = -83333972
$ws.Range("N99").Value = -5452.2
$ws.Range("H107").Value = 28446.4
$ws.Range("I107").Value = 2358
$ws.Range("K107").Value = 2358
$ws.Range("M107").Value = -438

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1718.5306
$ws.Range("I31").Value = 736.875
$ws.Range("J31").Value = 2194.4849
$ws.Range("K31").Value = 736.875
$ws.Range("L31").Value = 2194.4849
$ws.Range("M31").Value = -441.875
$ws.Range("N31").Value = -2784.4849
$ws.Range("H34").Value = 1718.5306
$ws.Range("I34").Value = 736.875
$ws.Range("J34").Value = 2194.4849
$ws.Range("K34").Value = 736.875
$ws.Range("L34").Value = 2194.4849
$ws.Range("M34").Value = -534.875
$ws.Range("N34").Value = -2598.4849
$ws.Range("H107").Value = 373.3793
$ws.Range("I107").Value = 326.8889
$ws.Range("K107").Value = 326.8889
$ws.Range("M107").Value = 1593.1111
$ws.Range("H140").Value = 47034.547
$ws.Range("J140").Value = 47034.547
$ws.Range("L140").Value = 47034.547
$ws.Range("N140").Value = -57394.547

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1221.35
$ws.Range("I68").Value = 647.1539
$ws.Range("J68").Value = 1588.459
$ws.Range("K68").Value = 1941.4617
$ws.Range("L68").Value = 4765.377
$ws.Range("M68").Value = -1130.4617
$ws.Range("N68").Value = -6387.377
$ws.Range("H71").Value = 1221.35
$ws.Range("I71").Value = 647.1539
$ws.Range("J71").Value = 1588.459
$ws.Range("K71").Value = 5824.3851
$ws.Range("L71").Value = 14296.131
$ws.Range("M71").Value = -1768.3851
$ws.Range("N71").Value = -22408.131
$ws.Range("H75").Value = 565
$ws.Range("I75").Value = 586.6667
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 1760.0001
$ws.Range("L75").Value = 1500
$ws.Range("M75").Value = -762.0001
$ws.Range("N75").Value = -3496
$ws.Range("H78").Value = 565
$ws.Range("I78").Value = 586.6667
$ws.Range("J78").Value = 500
$ws.Range("K78").Value = 5280.0003
$ws.Range("L78").Value = 4500
$ws.Range("M78").Value = -288.0002999999997
$ws.Range("N78").Value = -14484
$ws.Range("H80").Value = 8684.25
$ws.Range("I80").Value = 7849
$ws.Range("J80").Value = 8962.666999999999
$ws.Range("K80").Value = 23547
$ws.Range("L80").Value = 26888.001
$ws.Range("M80").Value = -22611
$ws.Range("N80").Value = -28760.001
$ws.Range("H83").Value = 8684.25
$ws.Range("I83").Value = 7849
$ws.Range("J83").Value = 8962.666999999999
$ws.Range("K83").Value = 70641
$ws.Range("L83").Value = 80664.003
$ws.Range("M83").Value = -65961
$ws.Range("N83").Value = -90024.003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2357.9285
$ws.Range("I97").Value = 2625
$ws.Range("J97").Value = 2001.8334
$ws.Range("K97").Value = 2625
$ws.Range("L97").Value = 2001.8334
$ws.Range("M97").Value = -2129
$ws.Range("N97").Value = -2993.8334
$ws.Range("H107").Value = 540.5
$ws.Range("J107").Value = 549.5
$ws.Range("L107").Value = 549.5
$ws.Range("N107").Value = -4389.5
$ws.Range("H113").Value = 1346.2
$ws.Range("I113").Value = 1357.75
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1357.75
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 812.25
$ws.Range("N113").Value = -5640
$ws.Range("H122").Value = 13157894
$ws.Range("I122").Value = 13157894
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 39473682
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -39471232
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2777.7222
$ws.Range("I132").Value = 1427.8096
$ws.Range("K132").Value = 4283.4288
$ws.Range("M132").Value = -1753.4288
$ws.Range("H138").Value = 30009.084
$ws.Range("J138").Value = 30009.084
$ws.Range("L138").Value = 30009.084
$ws.Range("N138").Value = -40289.084

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9613
$ws.Range("I122").Value = 13542.571
$ws.Range("J122").Value = 2736.25
$ws.Range("K122").Value = 40627.713
$ws.Range("L122").Value = 8208.75
$ws.Range("M122").Value = -38177.713
$ws.Range("N122").Value = -13108.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1741.4
$ws.Range("I122").Value = 1551.5
$ws.Range("J122").Value = 2026.25
$ws.Range("K122").Value = 4654.5
$ws.Range("L122").Value = 6078.75
$ws.Range("M122").Value = -2204.5
$ws.Range("N122").Value = -10978.75
